# Generate Report for Handback
# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# values for the first data row (row 2) on both the "zh-cn" and "de-de"
# language sheets, reflecting a fresh report generation run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-11 20:46:22"
$wsZhCn.Range("H2").Value = "2016-03-11 20:46:46"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-11 20:46:25"
$wsDeDe.Range("H2").Value = "2016-03-11 20:46:51"
